$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits right after
#    the "Fecha:08/10/24" run.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Retitle the (old) heading
#    "Análisis de la Realidad Sociocultural y la Tecnología en la
#     Provincia de Formosa"
#    to the new heading
#    "Acceso y uso de la tecnología de la comunicación en la
#     provincia de Formosa"
# ------------------------------------------------------------------
$old = "Análisis de la Realidad Sociocultural y la Tecnología en la Provincia de Formosa"
$new = "Acceso y uso de la tecnología de la comunicación en la provincia de Formosa"

$find = $d.Content
$found = $find.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ------------------------------------------------------------------
# 3) Re-locate the range of the freshly written heading text and split
#    it into three runs, the way Word does while text is typed/edited
#    in separate steps, and drop the "_GoBack" bookmark back in between
#    "...provincia de F" and "ormosa".
# ------------------------------------------------------------------
if ($found) {
    $headingRange = $d.Content
    $headingRange.Find.Execute($new) | Out-Null
    $start = $headingRange.Start

    # A temporary bookmark placed between "A" and "cceso..." forces Word to
    # keep those two stretches of identical-formatted text as separate runs
    # (mirroring the run break visible in the authored document) without
    # leaving any bookmark markers behind once it is removed again.
    $splitPos = $d.Range($start + 1, $start + 1)
    $d.Bookmarks.Add("ZZTempSplit", $splitPos)

    # Put "_GoBack" back between "...provincia de F" and "ormosa".
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $goBackPos = $d.Range($start + 69, $start + 69)
    $d.Bookmarks.Add("_GoBack", $goBackPos)

    if ($d.Bookmarks.Exists("ZZTempSplit")) {
        $d.Bookmarks.Item("ZZTempSplit").Delete()
    }
}
